$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; everything from row 71 down shifts to row+1.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with a new weekly price record (copy of the
# original row 71's record with the updated week's figures).
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44790
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112021
$ws.Range("G71").Value = "Ají"
$ws.Range("H71").Value = "Inferno"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 45
$ws.Range("K71").Value = 17000
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = 17556
$ws.Range("N71").Value = "$/caja 12 kilos"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 1463
$ws.Range("Q71").Value = 12
$ws.Range("R71").Value = "Hortaliza"
